$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold/bordered/centered) onto the two new
# header cells so they reuse the same style index instead of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-45
$data = @(
    @(8,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(5,5),
    @(5,6),
    @(9,9),
    @(5,5),
    @(7,7),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(7,8),
    @(6,8),
    @(8,8),
    @(8,9),
    @(7,7),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(9,9),
    @(6,6),
    @(6,7),
    @(8,8),
    @(8,8),
    @(7,8),
    @(7,7),
    @(6,6),
    @(6,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(7,7),
    @(9,9),
    @(5,5),
    @(5,5),
    @(4,4),
    @(4,4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
